$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Bang_Phan_Cong_Nhiem_Vu")
$ws2.Range("B8").Value = 22110337
$ws2.Range("C8").Value = "Nguyễn Lý Hùng"
$ws2.Range("D8").Value = "Viết API lấy các category theo chiều ngang , viết logic gọi api hiển thị thông tin category lên màn hình, merge code"
$ws2.Range("B9").Value = 22110369
$ws2.Range("C9").Value = "Lê Đình Lộc"
$ws2.Range("D9").Value = "Viết API lấy thông tin user , gọi api hiển thị thông tin user lên màn hình"
$ws2.Range("B10").Value = 22110379
$ws2.Range("C10").Value = "Võ Văn Nam"
$ws2.Range("D10").Value = "Thiết kế giao diện  màn hình chính, tinh chỉnh bottomAppBar"
$ws2.Range("B11").Value = 22110418
$ws2.Range("C11").Value = "Nguyễn Tuấn Thành"
$ws2.Range("D11").Value = "Xây dựng chức năng lấy tất cả sản phẩm theo từng category được sắp xếp tăng dần theo giá bán hiển thị dạng lưới"
$ws2.Range("B12").Value = 22110436
$ws2.Range("C12").Value = "Huỳnh Thái Toàn"
$ws2.Range("D12").Value = "Thiết kế giao diện đăng nhập và viết API đăng nhập"
$ws2.Range("D5").Value = "Thiết kế giao diện Intro, giao diện bottomAppBar, logic nút bắt đầu, login activity, merge code, tạo repo"